$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (price & volume refresh, plus a few
# coins that moved rank positions in the source ranking).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.289.05'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.241.37'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.80'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.80%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.85%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.72'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0810'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.17'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.53%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.384.10'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.25%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.582.64'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.829'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.50'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '44.008.37'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0963'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.39'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.12'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -7.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.50'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '237.70'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.94'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '38.55'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +5.62%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.21'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.02%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.90'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.01'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.83'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '153.11'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0793'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.28%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.13'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.53%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.82%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.76'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.50'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.80'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.27'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -7.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0299'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.71%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.745.58'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.48%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '82.91'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.33%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '14.97'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +9.08%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '99.63'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.92'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -4.81%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.07'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.64'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.97%  '
